# Update "想去人数" (want-to-go count) figures in column F across the
# four worksheets of the workbook, matching the refreshed data snapshot.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: 展览 ---------------------------------------------------
$ws = $wb.Worksheets.Item(1)
$ws.Range("F3").Value  = 168
$ws.Range("F5").Value  = 3346
$ws.Range("F6").Value  = 1130
$ws.Range("F7").Value  = 2248
$ws.Range("F9").Value  = 1136
$ws.Range("F12").Value = 1698
$ws.Range("F17").Value = 247
$ws.Range("F18").Value = 1609
$ws.Range("F19").Value = 268
$ws.Range("F20").Value = 1326
$ws.Range("F22").Value = 272
$ws.Range("F23").Value = 632
$ws.Range("F24").Value = 12390
$ws.Range("F25").Value = 12437
$ws.Range("F30").Value = 36
$ws.Range("F31").Value = 409
$ws.Range("F32").Value = 1936
$ws.Range("F33").Value = 4
$ws.Range("F35").Value = 213
$ws.Range("F36").Value = 622

# --- Sheet 2: 演出 -----------------------------------------------------
$ws = $wb.Worksheets.Item(2)
$ws.Range("F10").Value = 52

# --- Sheet 3: 本地生活 -------------------------------------------------
$ws = $wb.Worksheets.Item(3)
$ws.Range("F3").Value = 104

# --- Sheet 4: 全部类型 (combined view of all entries) ------------------
$ws = $wb.Worksheets.Item(4)
$ws.Range("F4").Value  = 168
$ws.Range("F6").Value  = 3346
$ws.Range("F7").Value  = 1130
$ws.Range("F8").Value  = 2248
$ws.Range("F10").Value = 1136
$ws.Range("F12").Value = 104
$ws.Range("F14").Value = 1698
$ws.Range("F22").Value = 247
$ws.Range("F23").Value = 1609
$ws.Range("F24").Value = 268
$ws.Range("F25").Value = 1326
$ws.Range("F27").Value = 272
$ws.Range("F29").Value = 632
$ws.Range("F30").Value = 12391
$ws.Range("F31").Value = 12437
$ws.Range("F36").Value = 36
$ws.Range("F37").Value = 409
$ws.Range("F40").Value = 1936
$ws.Range("F41").Value = 4
$ws.Range("F45").Value = 213
$ws.Range("F46").Value = 622
$ws.Range("F47").Value = 52
